$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph that follows the title heading.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2) Insert a new bold paragraph "Play Aloha! Cluster Pays Free - Review & Rating"
#    right before the final (italic) image-prompt paragraph, using raw OOXML so the
#    run/paragraph structure matches exactly (leading empty run + bold run).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newHeadingXml = '<w:p ' + $xmlNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aloha! Cluster Pays Free - Review &amp; Rating</w:t></w:r></w:p>'
$paddingXml = '<w:p ' + $xmlNs + '/>'
$insertionPoint.InsertXML($newHeadingXml + $paddingXml)

# InsertXML merges its final inserted paragraph mark with what followed, so the
# paragraph right after our new bold paragraph is a leftover empty paragraph -
# remove it to get back to the original (untouched) final paragraph.
$junkPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$junkPara.Range.Delete()

# 3) Replace the image-generation-prompt text of the final paragraph with the
#    meta description text (keeps the paragraph's existing italic run formatting).
$oldText = "Create a fun feature image of a happy Maya warrior wearing glasses, surrounded by colorful flowers and palm trees against a backdrop of a white sandy beach and a distant volcano. The image should be in cartoon style and depict the vibrant and lively theme of Aloha! Cluster Pays. The Maya warrior could be holding a Tiki sculpture, and the pineapples, seashells, and other colorful symbols from the game could be scattered around the scene. The aim is to showcase the upbeat and unusual gameplay experience of Aloha! Cluster Pays and highlight the tropical island theme of the game."
$newText = "Full review of Aloha! Cluster Pays online slot game. Learn about gameplay, unique features, and payouts. Play free online."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
